$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (Total) summary sheet: insert a new row for 2022-Q4
#    above the existing 2022-Q3 row, shifting everything else down by one
#    and renumbering the sequential index in column A.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Copy formatting from the (now shifted) row below so the new row matches
# the existing look (borders/fonts/alignment already used on this sheet).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.1

# Renumber the sequential index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q4" worksheet, positioned right after "总计".
#    Duplicate the "2022-Q3" sheet (so formatting/styles match exactly) and
#    then trim it down to a single fund row with the new Q4 figures.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Drop the extra fund rows (2022-Q3 had 5, Q4 only has 1).
$q4.Range("A3:H6").EntireRow.Delete()

# Row 2 (000763 / 工银新财富灵活配置混合) keeps its code & name; only the
# numeric-looking text columns and the rank change. Use a leading apostrophe
# so these remain text cells (matching the source data's string formatting)
# instead of being coerced into real numbers.
$q4.Range("D2").Value = "'2.82"
$q4.Range("E2").Value = "'93.80"
$q4.Range("F2").Value = "'3.55"
$q4.Range("G2").Value = "'0.1001"
$q4.Range("H2").Value = 7

# ---------------------------------------------------------------------------
# 3) Keep "2022-Q1" as the active/selected tab, matching the source file.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
